$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.951.25'
$ws.Range("E2").Value = '  +0.78%  '

$ws.Range("D3").Value = '1.813.23'
$ws.Range("E3").Value = '  +1.67%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '310.51'
$ws.Range("E5").Value = '  +0.15%  '

$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("D7").Value = '0.5002'
$ws.Range("E7").Value = '  -2.13%  '

$ws.Range("D8").Value = '0.3915'
$ws.Range("E8").Value = '  +1.18%  '

$ws.Range("D9").Value = '0.09730'
$ws.Range("E9").Value = '  +24.59%  '

$ws.Range("D10").Value = '1.100'
$ws.Range("E10").Value = '  +1.30%  '

$ws.Range("D11").Value = '40.87'

$ws.Range("D12").Value = '6.403'
$ws.Range("E12").Value = '  +3.20%  '

$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").Value = '20.45'
$ws.Range("E13").Value = '  +1.48%  '

$ws.Range("B14").Value = 'BinanceUSD'
$ws.Range("C14").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D14").Value = '1.001'
$ws.Range("E14").Value = '  -0.06%  '

$ws.Range("D15").Value = '1.815.68'
$ws.Range("E15").Value = '  +2.18%  '

$ws.Range("D16").Value = '7.264'
$ws.Range("E16").Value = '  +1.02%  '

$ws.Range("D17").Value = '0.00001138'
$ws.Range("E17").Value = '  +5.87%  '

$ws.Range("D18").Value = '92.40'
$ws.Range("E18").Value = '  +1.21%  '

$ws.Range("D19").Value = '0.06652'
$ws.Range("E19").Value = '  +1.58%  '

$ws.Range("D21").Value = '17.17'
$ws.Range("E21").Value = '  +1.05%  '

$ws.Range("D22").Value = '5.916'
$ws.Range("E22").Value = '  +0.30%  '

$ws.Range("D23").Value = '28.004.82'
$ws.Range("E23").Value = '  +0.79%  '

$ws.Range("D24").Value = '11.07'
$ws.Range("E24").Value = '  +0.88%  '

$ws.Range("D25").Value = '2.247'
$ws.Range("E25").Value = '  +0.97%  '

$ws.Range("D26").Value = '158.62'
$ws.Range("E26").Value = '  -1.13%  '

$ws.Range("D27").Value = '2.020.74'
$ws.Range("E27").Value = '  +1.77%  '

$ws.Range("D28").Value = '20.56'
$ws.Range("E28").Value = '  +2.12%  '

$ws.Range("D29").Value = '2.392'
$ws.Range("E29").Value = '  +1.38%  '

$ws.Range("D30").Value = '126.75'
$ws.Range("E30").Value = '  +2.72%  '

$ws.Range("D31").Value = '0.1064'
$ws.Range("E31").Value = '  -1.19%  '

$ws.Range("D32").Value = '1.032'
$ws.Range("E32").Value = '  +0.08%  '

$ws.Range("D33").Value = '5.558'
$ws.Range("E33").Value = '  +1.62%  '

$ws.Range("D34").Value = '3.583'
$ws.Range("E34").Value = '  -1.35%  '

$ws.Range("D35").Value = '0.06735'
$ws.Range("E35").Value = '  -4.23%  '

$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").Value = '8.912'
$ws.Range("E36").Value = '  +1.91%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.02325'
$ws.Range("E37").Value = '  +1.02%  '

$ws.Range("D38").Value = '0.2141'
$ws.Range("E38").Value = '  +1.01%  '

$ws.Range("D39").Value = '4.938'
$ws.Range("E39").Value = '  -0.78%  '

$ws.Range("D40").Value = '11.26'
$ws.Range("E40").Value = '  -1.50%  '

$ws.Range("D41").Value = '0.6174'
$ws.Range("E41").Value = '  +1.67%  '

$ws.Range("D42").Value = '1.168'
$ws.Range("E42").Value = '  +2.03%  '

$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  -0.08%  '

$ws.Range("D44").Value = '13.20'
$ws.Range("E44").Value = '  +1.13%  '

$ws.Range("D45").Value = '0.5890'
$ws.Range("E45").Value = '  +0.28%  '

$ws.Range("B46").Value = 'WEMIXTOKEN'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '1.285'
$ws.Range("E46").Value = '  -2.65%  '

$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = '3.690'
$ws.Range("E47").Value = '  -0.15%  '

$ws.Range("D48").Value = '123.73'
$ws.Range("E48").Value = '  -0.60%  '

$ws.Range("D49").Value = '1.936'
$ws.Range("E49").Value = '  +2.10%  '

$ws.Range("D50").Value = '1.178'
$ws.Range("E50").Value = '  -1.84%  '

$ws.Range("D51").Value = '0.06771'
$ws.Range("E51").Value = '  -0.95%  '

